# Add a new course row ("Understanding Data Visualization") to the
# datacamp tracking sheet, just below the existing "Understanding Data
# Science" / "Understanding Data Engineering" rows (rows 41-42), pushing
# the remaining course rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 43 (Excel shifts rows 43.. down to 44.. and
# inherits formatting from the surrounding rows, same as a manual
# right-click > Insert in the UI).
$ws.Rows("43:43").Insert()

# Fill in the new course's data: course name in column A, rating in
# column I (matching the layout used by the two rows directly above it).
$ws.Range("A43").Value = "Understanding Data Visualization"
$ws.Range("I43").Value = 3

# Leave the selection where the user ended up after inserting the row.
$ws.Range("A49").Select()
